$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Sheet4")

# New simulation values (Knapsack bugfix reran the GA -> fewer generations, new fitness numbers)
$bVals = @(261,261,261,261,261,261,261,261,261,261,261,261,261,261,261,261,261,261,261)
$cVals = @(115.54,130.72,152.21,161.53,215.06,250.9,242.65,237.25,248.14,241.28,239.26,243.29,238.29,246.07,241.51,239.62,247.86,244.02,244.65)

for ($i = 0; $i -lt $bVals.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 2).Value = $bVals[$i]
    $ws4.Cells.Item($row, 3).Value = $cVals[$i]
}

# Row 21 no longer holds a generation (the run converged one generation earlier)
$ws4.Range("A21:C21").ClearContents()

# Apply the number format the author picked when re-checking the Generation column
$ws4.Range("A3:A21").NumberFormat = "#,##0.00"

# Make Sheet4 the active sheet/tab and select the data block, matching the reviewed range
$ws4.Activate()
$ws4.Range("A2:C20").Select()
